# Uber UK operations.docx -- add more cities, split two date cells into
# a "short month" run + a " year" run.
#
# Strategy: Range.InsertXML() lets us drop exact WordprocessingML
# fragments in place (replacing the target Range's contents), so we can
# build precisely the <w:p>/<w:r> shapes the diff calls for instead of
# relying on Find/Replace (keeps one run) or naive InsertBefore/Text
# edits (which pick up stray <w:rPr> overrides when a run is split via
# formatting toggles).

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1) "July 2012"    -> "Jul" + " 2012"   (two runs, same paragraph)
# 2) "December 2015" -> "Dec" + " 2015"  (two runs, same paragraph)
# ---------------------------------------------------------------------
function Split-DateCell($fullText, $shortMonth, $yearPart) {
    $hit = $d.Content
    $hit.Find.Execute($fullText, $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0)
    $xml = "<w:p $wNs><w:r><w:t>$shortMonth</w:t></w:r>" + `
           "<w:r><w:t xml:space=`"preserve`"> $yearPart</w:t></w:r></w:p>"
    $hit.InsertXML($xml)
}

Split-DateCell "July 2012" "Jul" "2012"
Split-DateCell "December 2015" "Dec" "2015"

# ---------------------------------------------------------------------
# 3) Add the new city paragraphs after the table, right before the
#    trailing bookmark ("_GoBack") paragraph that closes the document.
# ---------------------------------------------------------------------
function P([string]$innerXml = "") {
    return "<w:p $wNs>$innerXml</w:p>"
}

function TextRun([string]$city, [bool]$pageBreak = $false) {
    $brk = ""
    if ($pageBreak) { $brk = "<w:lastRenderedPageBreak/>" }
    return P "<w:r>$brk<w:t>$city</w:t></w:r>"
}

function CenteredEmpty() {
    return P "<w:pPr><w:jc w:val=`"center`"/></w:pPr>"
}

$firstBatch  = @("Bradford","Liverpool","Bristol","Wakefield","Coventry","Sunderland","Brighton","Newcastle")
$secondBatch = @("Liverpool","Bristol","Wakefield","Coventry","Sunderland","Brighton","Newcastle")

$pieces = New-Object System.Collections.ArrayList

[void]$pieces.Add((P))                                   # leading blank paragraph
foreach ($city in $firstBatch) { [void]$pieces.Add((TextRun $city)) }

for ($i = 0; $i -lt 9; $i++) { [void]$pieces.Add((P)) }   # nine blank paragraphs
[void]$pieces.Add((CenteredEmpty))                        # two centered-empty paragraphs
[void]$pieces.Add((CenteredEmpty))

[void]$pieces.Add((TextRun $secondBatch[0] $true))        # "Liverpool" w/ lastRenderedPageBreak
for ($i = 1; $i -lt $secondBatch.Count; $i++) {
    [void]$pieces.Add((TextRun $secondBatch[$i]))
}

# InsertXML("<w:p>..</w:p><w:p>..</w:p>") on a collapsed Range sitting
# right at the start of an existing paragraph merges the *last*
# fragment's runs into that existing paragraph instead of keeping it
# separate. Tack on one extra empty <w:p> so that merge "eats" a blank
# paragraph instead of one we actually care about, leaving the real
# _GoBack paragraph completely untouched.
[void]$pieces.Add((P))

$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertPoint.InsertXML(($pieces -join ""))
